# Redeem points 71277620 76.0
#
# The redemptions log's last row (A16, previously the phone number stored
# as text) is normalized to a real number, and a brand-new redemption row
# (17) is appended for phone 71277620 redeeming 76 points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: the phone number that was stored as text becomes a plain number.
$ws.Range("A16").Value = 71277620

# Row 17 (new): phone number stays text (matches how new rows are logged),
# points is numeric, timestamp is text.
$ws.Range("A17").Value = "'71277620"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = 76
$ws.Range("C17").Value = "2025-08-18T17:04:40"
